$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aguilar Prototype")

# Row 29: T5YIFR
$ws.Range("N29").Value = 45989
$ws.Range("R29").Value = $null
$ws.Range("S29").Value = 2.17
$ws.Range("T29").Value = 2.17
$ws.Range("U29").Value = 2.16

# Row 30: T10YIE
$ws.Range("N30").Value = 45989
$ws.Range("R30").Value = $null
$ws.Range("T30").Value = 2.22
$ws.Range("U30").Value = 2.23

# Row 47: DFF
$ws.Range("N47").Value = 45988

# Row 48: DGS2
$ws.Range("N48").Value = 45987
$ws.Range("Q48").Value = 3.45
$ws.Range("R48").Value = 3.43
$ws.Range("S48").Value = 3.46
$ws.Range("U48").Value = $null

# Row 49: DGS5
$ws.Range("N49").Value = 45987
$ws.Range("Q49").Value = 3.56
$ws.Range("R49").Value = 3.55
$ws.Range("S49").Value = 3.61
$ws.Range("U49").Value = $null

# Row 50: DGS10
$ws.Range("N50").Value = 45987
$ws.Range("Q50").Value = 4
$ws.Range("R50").Value = 4.01
$ws.Range("S50").Value = 4.04
$ws.Range("U50").Value = $null

# Row 52: DBAA
$ws.Range("N52").Value = 45987
$ws.Range("Q52").Value = 5.78
$ws.Range("R52").Value = 5.8
$ws.Range("S52").Value = 5.84
$ws.Range("U52").Value = $null
